$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Un-merge the three header blocks on row 47 (K47:P47 / R47:W47 / AA47:AF47)
#    The target layout no longer merges these ranges at all.
# ---------------------------------------------------------------------------
$ws.Range("K47:P47").UnMerge()
$ws.Range("R47:W47").UnMerge()
$ws.Range("AA47:AF47").UnMerge()

# ---------------------------------------------------------------------------
# 2) Row 47 : K47/L47, R47/S47, AA47/AB47
# ---------------------------------------------------------------------------
# Build the new "closing cell, no alignment" style (will become cellXfs index 68)
# by grabbing the existing right-border style from P47 and stripping its
# centered alignment.
$ws.Range("P47").Copy()
$ws.Range("L47").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("L47").HorizontalAlignment = 1      # xlGeneral -> drop the centering
$ws.Range("L47").ClearContents()

# K47 keeps its text (shared string 40) but switches to the "no alignment"
# left-border style already used by D48.
$ws.Range("D48").Copy()
$ws.Range("K47").PasteSpecial(-4122)

# Propagate the same two styles to the R/S and AA/AB blocks.
$ws.Range("L47").Copy()
$ws.Range("S47").PasteSpecial(-4122)
$ws.Range("AB47").PasteSpecial(-4122)

$ws.Range("K47").Copy()
$ws.Range("R47").PasteSpecial(-4122)
$ws.Range("AA47").PasteSpecial(-4122)

# Drop the now unused trailing cells of each block.
$ws.Range("M47:P47").Clear()
$ws.Range("T47:W47").Clear()
$ws.Range("AC47:AF47").Clear()

# ---------------------------------------------------------------------------
# 3) Row 48 : L48/S48/AB48 adopt the "no alignment" right border style
#    (the one already used at H32), keeping their existing text.
# ---------------------------------------------------------------------------
$ws.Range("H32").Copy()
$ws.Range("L48").PasteSpecial(-4122)
$ws.Range("S48").PasteSpecial(-4122)
$ws.Range("AB48").PasteSpecial(-4122)

$ws.Range("M48:P48").Clear()
$ws.Range("T48:W48").Clear()
$ws.Range("AC48:AF48").Clear()

# ---------------------------------------------------------------------------
# 4) Row 49
# ---------------------------------------------------------------------------
# K49 / R49 / AA49 keep their style but lose their numeric value.
$ws.Range("K49").ClearContents()
$ws.Range("R49").ClearContents()
$ws.Range("AA49").ClearContents()

# L49 / S49 take on the style already used at AF37, keeping their value (4).
$ws.Range("AF37").Copy()
$ws.Range("L49").PasteSpecial(-4122)
$ws.Range("S49").PasteSpecial(-4122)

# AB49 takes on the style already used at AF34, keeping its value (2).
$ws.Range("AF34").Copy()
$ws.Range("AB49").PasteSpecial(-4122)

$ws.Range("M49:P49").Clear()
$ws.Range("T49:W49").Clear()
$ws.Range("AC49:AF49").Clear()

# ---------------------------------------------------------------------------
# 5) Row 50 : K50/L50 disappear entirely, AA50 loses its value, AB50 takes the
#    "H33" style while keeping its value (4).
# ---------------------------------------------------------------------------
$ws.Range("K50:L50").Clear()
$ws.Range("AA50").ClearContents()
$ws.Range("H33").Copy()
$ws.Range("AB50").PasteSpecial(-4122)
$ws.Range("AC50:AF50").Clear()

# ---------------------------------------------------------------------------
# 6) Row 51 : same shape as row 50 (value 5).
# ---------------------------------------------------------------------------
$ws.Range("K51:L51").Clear()
$ws.Range("AA51").ClearContents()
$ws.Range("H33").Copy()
$ws.Range("AB51").PasteSpecial(-4122)
$ws.Range("AC51:AF51").Clear()

# ---------------------------------------------------------------------------
# 7) Row 52 : AA52 loses its value, AB52 takes the AF37 style (value 6).
# ---------------------------------------------------------------------------
$ws.Range("AA52").ClearContents()
$ws.Range("AF37").Copy()
$ws.Range("AB52").PasteSpecial(-4122)
$ws.Range("AC52:AF52").Clear()

# ---------------------------------------------------------------------------
# 8) New column widths for S (19) and AB (28).
# ---------------------------------------------------------------------------
$ws.Columns("S").ColumnWidth = 12.666666666666666
$ws.Columns("AB").ColumnWidth = 13.5

# ---------------------------------------------------------------------------
# 9) Selection / active cell as last left by the author.
# ---------------------------------------------------------------------------
$ws.Range("AA56").Select()
